# "Undid the italic on the header."
#
# The document opens with four centered, bold ("Strong" style), blue lines
# that together form the header block:
#   1. Marlon Torres
#   2. 11/27/2013
#   3. Web Programming Fundamentals - Section 01
#   4. Activity: Problem Solving
#
# Every one of those four lines currently carries explicit italic (w:i) on
# the run that holds the visible text, and the first three also carry it on
# the paragraph mark's own run properties (pPr/rPr) - the "Activity: Problem
# Solving" line's paragraph mark was never italicized, only its text run
# was. We turn italics off everywhere it is actually set, without touching
# paragraph marks that never had it (so we don't invent new formatting that
# wasn't part of the original document / the requested change).

$d = $word.ActiveDocument

$headerLines = @(
    "Marlon Torres",
    "11/27/2013",
    "Web Programming Fundamentals",
    "Activity: Problem Solving"
)

foreach ($para in $d.Paragraphs) {
    $range = $para.Range
    $text = $range.Text

    $matched = $false
    foreach ($line in $headerLines) {
        if ($text.Contains($line)) {
            $matched = $true
        }
    }
    if (-not $matched) {
        continue
    }

    if ($text.Contains("Activity: Problem Solving")) {
        # Only the run (visible text) is italic here - leave the paragraph
        # mark's formatting exactly as it was.
        $textOnly = $d.Range($range.Start, $range.End - 1)
        $textOnly.Font.Italic = $false
    } else {
        # Both the run and the paragraph mark are italic - clear the whole
        # paragraph range (this covers pPr/rPr as well as the run's rPr).
        $range.Font.Italic = $false
    }
}
